$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.177.23"
$ws.Range("E2").Value = "  -1.25%  "

$ws.Range("D3").Value = "2.399.86"
$ws.Range("E3").Value = "  -2.08%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "505.63"
$ws.Range("E5").Value = "  -3.72%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.68"
$ws.Range("E6").Value = "  +1.47%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.24%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.558"
$ws.Range("E8").Value = "  -0.79%  "

$ws.Range("D9").Value = "2.433.99"
$ws.Range("E9").Value = "  -0.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0973"
$ws.Range("E10").Value = "  -0.65%  "

$ws.Range("E11").Value = "  -0.93%  "

$ws.Range("E12").Value = "  +0.21%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.66"
$ws.Range("E13").Value = "  -5.80%  "

$ws.Range("D14").Value = "2.845.59"
$ws.Range("E14").Value = "  -1.32%  "

$ws.Range("D15").Value = "57.079.87"
$ws.Range("E15").Value = "  -1.24%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.82"
$ws.Range("E16").Value = "  +0.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000133"
$ws.Range("E17").Value = "  -0.01%  "

$ws.Range("D18").Value = "2.435.41"
$ws.Range("E18").Value = "  -0.54%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.29"
$ws.Range("E19").Value = "  -1.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.11"
$ws.Range("E20").Value = "  -1.11%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "312.73"
$ws.Range("E21").Value = "  -1.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.53"
$ws.Range("E22").Value = "  +7.16%  "

$ws.Range("E23").Value = "  -0.30%  "

$ws.Range("E24").Value = "  -1.83%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.98"
$ws.Range("E25").Value = "  -0.13%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.995"
$ws.Range("E26").Value = "  -0.51%  "

$ws.Range("D27").Value = "2.504.27"
$ws.Range("E27").Value = "  -2.28%  "

$ws.Range("E28").Value = "  -6.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.152"
$ws.Range("E29").Value = "  -2.93%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.46"
$ws.Range("E30").Value = "  +2.92%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.58"
$ws.Range("E31").Value = "  +0.20%  "

$ws.Range("D32").Value = "0.0₃0733"
$ws.Range("E32").Value = "  -0.68%  "

$ws.Range("E33").Value = "  -0.59%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.20"
$ws.Range("E34").Value = "  +1.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.13"
$ws.Range("E35").Value = "  -1.74%  "

$ws.Range("E36").Value = "  -0.14%  "

$ws.Range("E37").Value = "  -0.52%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.06"
$ws.Range("E38").Value = "  +1.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.22"
$ws.Range("E39").Value = "  +3.30%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.83"
$ws.Range("E40").Value = "  +0.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.818"
$ws.Range("E41").Value = "  +2.75%  "

$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "36.18"
$ws.Range("E42").Value = "  -0.28%  "

$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.45"
$ws.Range("E43").Value = "  -0.97%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "132.04"
$ws.Range("E44").Value = "  +5.97%  "

$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.40"
$ws.Range("E45").Value = "  -0.61%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.00"
$ws.Range("E46").Value = "  +3.68%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "257.75"
$ws.Range("E47").Value = "  -2.87%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.568"
$ws.Range("E48").Value = "  -2.72%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0913"
$ws.Range("E49").Value = "  -1.78%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0496"
$ws.Range("E50").Value = "  +0.12%  "

$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0213"
$ws.Range("E51").Value = "  +0.60%  "
